$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$edits = @(
  @("B2", "건축학부"),
  @("C2", "29"),
  @("D2", "9"),
  @("E2", $null),
  @("B3", "건축공학부"),
  @("D3", "2"),
  @("E3", $null),
  @("B4", "건설환경공학과"),
  @("D4", "5"),
  @("E4", $null),
  @("B5", "도시공학과"),
  @("D5", "3"),
  @("E5", $null),
  @("B6", "자원환경공학과"),
  @("D6", "3"),
  @("E6", $null),
  @("B7", "융합전자공학부"),
  @("D7", "27"),
  @("E7", $null),
  @("B8", "컴퓨터소프트웨어학부"),
  @("D8", "31"),
  @("E8", $null),
  @("B9", "전기공학전공"),
  @("E9", $null),
  @("B10", "바이오메디컬공학전공"),
  @("D10", "10"),
  @("E10", $null),
  @("B11", "신소재공학부"),
  @("D11", "21"),
  @("E11", $null),
  @("B12", "화학공학과"),
  @("D12", "19"),
  @("E12", $null),
  @("B13", "생명공학과"),
  @("D13", "10"),
  @("E13", $null),
  @("B14", "유기나노공학과"),
  @("D14", "5"),
  @("E14", $null),
  @("B15", "에너지공학과"),
  @("D15", "5"),
  @("E15", $null),
  @("B16", "기계공학부"),
  @("D16", "30"),
  @("E16", $null),
  @("B17", "원자력공학과"),
  @("D17", "3"),
  @("E17", $null),
  @("B18", "산업공학과"),
  @("D18", "7"),
  @("E18", $null),
  @("B19", "미래자동차공학과"),
  @("D19", "6"),
  @("E19", $null),
  @("B20", "데이터사이언스학부"),
  @("D20", "12"),
  @("E20", $null),
  @("B21", "수학과"),
  @("D21", "7"),
  @("E21", $null),
  @("B22", "물리학과"),
  @("D22", "10"),
  @("E22", $null),
  @("B23", "화학과"),
  @("D23", "21"),
  @("E23", $null),
  @("B24", "생명과학과"),
  @("D24", "22"),
  @("E24", $null),
  @("B25", "의류학과"),
  @("D25", "9"),
  @("E25", $null),
  @("B26", "식품영양학과"),
  @("D26", "7"),
  @("E26", $null),
  @("B27", "실내건축디자인학과"),
  @("D27", "7"),
  @("E27", $null),
  @("B28", "간호학과"),
  @("D28", "16"),
  @("E28", $null),
  @("B29", "정보시스템학과"),
  @("C29", "24"),
  @("D29", "16"),
  @("E29", $null),
  @("B30", "국어국문학과"),
  @("D30", "6"),
  @("E30", $null),
  @("B31", "중어중문학과"),
  @("D31", "6"),
  @("E31", $null),
  @("B32", "영어영문학과"),
  @("D32", "6"),
  @("E32", $null),
  @("B33", "독어독문학과"),
  @("C33", $null),
  @("D33", "2"),
  @("E33", $null),
  @("B34", "사학과"),
  @("D34", "7"),
  @("B35", "철학과"),
  @("D35", "3"),
  @("B36", "정치외교학과"),
  @("D36", "8"),
  @("B37", "사회학과"),
  @("D37", "6"),
  @("B38", "미디어커뮤니케이션학과"),
  @("D38", "22"),
  @("B39", "관광학부"),
  @("D39", "6"),
  @("B40", "정책학과"),
  @("D40", "16"),
  @("E40", $null),
  @("B41", "행정학과"),
  @("D41", "4"),
  @("E41", $null),
  @("B42", "경제금융학부"),
  @("D42", "30"),
  @("E42", $null),
  @("B43", "경영학부"),
  @("D43", "51"),
  @("E43", $null),
  @("B44", "파이낸스경영학과"),
  @("D44", "14"),
  @("E44", $null),
  @("B45", "스포츠매니지먼트전공"),
  @("D45", "11"),
  @("E45", $null),
  @("A46", $null),
  @("B46", "연극영화학과(영화전공)"),
  @("C46", $null),
  @("D46", "9"),
  @("E46", $null),
  @("A47", $null),
  @("B47", "국제학전공"),
  @("C47", "3"),
  @("D47", "20")
)

foreach ($edit in $edits) {
  $ref = $edit[0]
  $val = $edit[1]
  if ($null -eq $val) {
    $ws.Range($ref).Value = ""
  } else {
    $ws.Range($ref).Value = $val
  }
}
